# Update template servizi czrm
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column C ("Email") to make room for the
# new "Orario al pubblico del telefono" header, then add a new trailing
# column "Link" after the current last column ("PEC").
$ws.Columns("C:C").Insert()
$ws.Range("F1").EntireColumn.Insert()

# Header row values.
$ws.Range("C1").Value = "Orario al pubblico del telefono"
$ws.Range("F1").Value = "Link"

# Match the style (yellow fill) used by the rest of the header row.
$ws.Range("C1").Interior.Color = $ws.Range("B1").Interior.Color
$ws.Range("F1").Interior.Color = $ws.Range("E1").Interior.Color

# Column C width per the target column definition (autofit to its header text).
$ws.Columns("C:C").ColumnWidth = 27.833333333333332

# Selection moves to F2.
$ws.Range("F2").Select()
